# Append two new India rows ("ladakh" and "Jammu and Kashmir") below the
# existing state list (last existing data row is 29 / "West Bengal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values.
$ws.Range("A30").Value = "India"
$ws.Range("B30").Value = "ladakh"
$ws.Range("A31").Value = "India"
$ws.Range("B31").Value = "Jammu and Kashmir"

# Match the formatting of the preceding data row (style index used by
# column A, which is what the new rows end up with in the source file)
# by copying A29's format onto the four new cells.
$ws.Range("A29").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("B30").PasteSpecial(-4122)
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("B31").PasteSpecial(-4122)
